# Agrega planificacion con plastificado
# Adds a new "Plastificado" row to the "Maquinas" sheet (as the new row 6,
# pushing the existing rows down) and makes "Maquinas" the active sheet,
# replacing "OrdenEstandar" as the previously active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maquinas")

# Remove the stray formatted row far below the data (row 1048576) before
# inserting, so it does not get shifted down and left behind as clutter.
$ws.Rows(1048576).Delete()

# Insert a new row at position 6 and populate it with the new machine info.
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = "Plastificado"
$ws.Range("B6").Value = "Plastificadora"
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10

# Make "Maquinas" the active sheet (this also clears the previous
# tabSelected/topLeftCell state on "OrdenEstandar") and set the selection.
$ws.Activate()
[void]$ws.Range("A7").Select()
